$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4046949148178101
$ws.Range("B1").Value = 0.366361141204834
$ws.Range("C1").Value = 3.402585029602051
$ws.Range("D1").Value = 1.604530334472656
$ws.Range("E1").Value = 1.131492972373962
